{"js": "// Locate the final acceptance line: \"I accept the WBA (Sravan Krsna Rao)\"\nconst results = context.document.body.search(\n  \"I accept the WBA (Sravan Krsna Rao)\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Replace with the \"revised\" wording, split into the same run layout Word\n  // produces when text is typed mid-sentence: the existing wording becomes\n  // two runs around the newly typed \"revised \", followed by the \"_GoBack\"\n  // bookmark Word drops at the last edited location.\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    \"<w:p>\" +\n    '<w:r><w:t xml:space=\"preserve\">I accept the </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">revised </w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>WBA (Sravan Krsna Rao)</w:t></w:r>\" +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  target.insertOoxml(ooxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the final acceptance line: \"I accept the WBA (Sravan Krsna Rao)\"\n$f = $d.Content\n$f.Find.ClearFormatting()\n$found = $f.Find.Execute(\"I accept the WBA (Sravan Krsna Rao)\")\n\nif ($found) {\n    # Build a plain Range (not the live Find range) over the matched text so\n    # InsertXML replaces exactly that span.\n    $target = $d.Range($f.Start, $f.End)\n\n    # Replace with the \"revised\" wording, split into the same run layout Word\n    # produces when text is typed mid-sentence: the existing wording becomes\n    # two runs around the newly typed \"revised \", followed by the \"_GoBack\"\n    # bookmark Word drops at the last edited location.\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n        '<w:r><w:t xml:space=\"preserve\">I accept the </w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">revised </w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>WBA (Sravan Krsna Rao)</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $target.InsertXML($xml)\n}\n"}
